$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 81
$ws.Range("A3").Value = 82
$ws.Range("A4").Value = 83
$ws.Range("A5").Value = 84
$ws.Range("A6").Value = 85
$ws.Range("A7").Value = 86
$ws.Range("A8").Value = 87
$ws.Range("A9").Value = 88
$ws.Range("A10").Value = 89
$ws.Range("A11").Value = 90
$ws.Range("A12").Value = 91
